$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "dadosDeAcesso"

# Header row (order matters for shared-string table allocation)
$ws.Range("A1").Value = "id"
$ws.Range("C1").Value = "senha"
$ws.Range("B1").Value = "nomeDeUsuário"

# Data row
$ws.Range("A2").Value = "ID_0001"
$ws.Range("B2").Value = "André Automatizador"
$ws.Range("C2").Value = "automacaoteste"

# Autofit columns B and C to match bestFit widths from the target workbook
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null

$ws.Range("C2").Select()
